$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 3 (shifts old row 3 "Unassigned" down to row 4)
$ws.Rows.Item(3).Insert()

# Row 2: change species name
$ws.Range("A2").Value = "Fundulus heteroclitus or majalis"

# Row 3: newly inserted row gets the old Row2 value
$ws.Range("A3").Value = "Cololabis saira"

# Row 4 already has "Unassigned" in A4/B4/C4 from the shifted row, ensure values are correct
$ws.Range("A4").Value = "Unassigned"
$ws.Range("B4").Value = "Unassigned"
$ws.Range("C4").Value = "Unassigned"

# New rows 5 and 6 with additional species
$ws.Range("A5").Value = "Mareca americana"
$ws.Range("A6").Value = "Myrophis vafer"
